$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.258.31'
$ws.Range('E2').Value = '  +0.62%  '
$ws.Range('D3').Value = '1.659.90'
$ws.Range('E3').Value = '  +0.55%  '
$ws.Range('E4').Value = '  +0.66%  '
$ws.Range('D5').Value = "'218.32"
$ws.Range('E5').Value = '  +0.04%  '
$ws.Range('D6').Value = "'0.5323"
$ws.Range('E6').Value = '  +0.57%  '
$ws.Range('E7').Value = '  +0.64%  '
$ws.Range('E8').Value = '  +0.96%  '
$ws.Range('D9').Value = "'0.06344"
$ws.Range('E9').Value = '  +0.62%  '
$ws.Range('E10').Value = '  +0.77%  '
$ws.Range('D11').Value = "'0.07833"
$ws.Range('E11').Value = '  +1.09%  '
$ws.Range('D12').Value = "'4.538"
$ws.Range('E12').Value = '  +1.58%  '
$ws.Range('D13').Value = '1.665.77'
$ws.Range('E13').Value = '  +0.88%  '
$ws.Range('D14').Value = '1.888.35'
$ws.Range('E14').Value = '  +0.60%  '
$ws.Range('D15').Value = "'0.5513"
$ws.Range('E15').Value = '  +1.20%  '
$ws.Range('D16').Value = '0.0₅8170'
$ws.Range('E16').Value = '  +0.73%  '
$ws.Range('E17').Value = '  +0.72%  '
$ws.Range('D18').Value = '26.251.51'
$ws.Range('E18').Value = '  +0.56%  '
$ws.Range('E19').Value = '  +0.70%  '
$ws.Range('D20').Value = "'4.647"
$ws.Range('E20').Value = '  +2.22%  '
$ws.Range('D21').Value = "'191.97"
$ws.Range('E21').Value = '  -0.65%  '
$ws.Range('E22').Value = '  +0.84%  '
$ws.Range('D23').Value = "'6.046"
$ws.Range('E23').Value = '  +1.23%  '
$ws.Range('E24').Value = '  +0.66%  '
$ws.Range('D25').Value = "'144.66"
$ws.Range('E25').Value = '  +3.30%  '
$ws.Range('D26').Value = "'0.1229"
$ws.Range('E26').Value = '  -0.82%  '
$ws.Range('D27').Value = "'7.225"
$ws.Range('E27').Value = '  -0.17%  '
$ws.Range('D28').Value = "'16.05"
$ws.Range('E28').Value = '  -0.49%  '
$ws.Range('D29').Value = "'1.471"
$ws.Range('E29').Value = '  +2.50%  '
$ws.Range('D30').Value = "'0.05787"
$ws.Range('E30').Value = '  -1.87%  '
$ws.Range('D31').Value = "'1.279"
$ws.Range('E31').Value = '  +0.06%  '
$ws.Range('D32').Value = "'3.569"
$ws.Range('E32').Value = '  +2.00%  '
$ws.Range('D33').Value = "'3.285"
$ws.Range('E33').Value = '  +1.68%  '
$ws.Range('E34').Value = '  +4.19%  '
$ws.Range('E35').Value = '  +2.24%  '
$ws.Range('D36').Value = "'0.9572"
$ws.Range('E36').Value = '  +1.56%  '
$ws.Range('E37').Value = '  +0.72%  '
$ws.Range('D38').Value = "'0.5791"
$ws.Range('E38').Value = '  +2.58%  '
$ws.Range('E39').Value = '  +0.05%  '
$ws.Range('D40').Value = "'5.842"
$ws.Range('E40').Value = '  +0.08%  '
$ws.Range('D41').Value = "'0.8508"
$ws.Range('E41').Value = '  +0.93%  '
$ws.Range('E42').Value = '  +0.65%  '
$ws.Range('E43').Value = '  +3.96%  '
$ws.Range('D44').Value = '1.044.91'
$ws.Range('E44').Value = '  +3.85%  '
$ws.Range('D45').Value = '1.801.54'
$ws.Range('E45').Value = '  +0.45%  '
$ws.Range('D46').Value = "'57.05"
$ws.Range('E46').Value = '  +0.55%  '
$ws.Range('D47').Value = '0.0₈105'
$ws.Range('E47').Value = '  -0.61%  '
$ws.Range('D48').Value = "'1.011"
$ws.Range('E48').Value = '  +0.56%  '
$ws.Range('E49').Value = '  +1.95%  '
$ws.Range('D50').Value = "'7.947"
$ws.Range('E50').Value = '  +2.20%  '
$ws.Range('D51').Value = "'0.05160"
$ws.Range('E51').Value = '  +0.19%  '
